# Appends a newly scraped Lancers.jp job listing as a new row (row 6) into the
# "ランサーズ" (Lancers) sheet, pushing the existing entries down by one row,
# and refreshes the "取得日時" (fetched-at) timestamp for every data row to the
# new run's timestamp (2025-11-25 01:20:23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-11-25 01:20:23"

# --- Insert a fresh row at position 6 (shifts old rows 6..18 down to 7..19) ---
$ws.Rows.Item(6).Insert()

# --- Refresh the "取得日時" timestamp for every data row (2..19) ---
for ($r = 2; $r -le 19; $r++) {
    $ws.Range("A$r").Value = $newTimestamp
}

# --- Populate the newly inserted row 6 with the new job listing ---
$ws.Range("B6").Value = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5405023"
$ws.Range("G6").Value = 178
$ws.Range("H6").Value = "★bot ◆ツール"

# --- Rebuild the URL hyperlinks for column F (F2:F19) from scratch.
#     (Row-insert carries the old hyperlink relationship onto the new row,
#     so the cleanest way to land on the correct final mapping - each F-cell
#     pointing at the right lancers.jp URL with a single relationship id - is
#     to clear every hyperlink and re-add them top-to-bottom.) ---
$ws.Hyperlinks.Delete()

$urls = @(
    "https://www.lancers.jp/work/detail/5434693",
    "https://www.lancers.jp/work/detail/5440400",
    "https://www.lancers.jp/work/detail/5423720",
    "https://www.lancers.jp/work/detail/5419380",
    "https://www.lancers.jp/work/detail/5405023",
    "https://www.lancers.jp/work/detail/5440052",
    "https://www.lancers.jp/work/detail/5440077",
    "https://www.lancers.jp/work/detail/5439484",
    "https://www.lancers.jp/work/detail/5431107",
    "https://www.lancers.jp/work/detail/5440318",
    "https://www.lancers.jp/work/detail/5440436",
    "https://www.lancers.jp/work/detail/5440417",
    "https://www.lancers.jp/work/detail/5440440",
    "https://www.lancers.jp/work/detail/5440230",
    "https://www.lancers.jp/work/detail/5440042",
    "https://www.lancers.jp/work/detail/5440002",
    "https://www.lancers.jp/work/detail/5440204",
    "https://www.lancers.jp/work/detail/5440325"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Range("F$row"), $urls[$i]) | Out-Null
    $ws.Range("F$row").Style = "Hyperlink"
}

Write-Host "Row insert + timestamp refresh + hyperlink rebuild complete"
